$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 250, pushing the existing rows 250..331
# down to 252..333 (preserving their formatting, e.g. the date style on
# column D).
$ws.Rows.Item(250).Resize(2).Insert()

# New row 250 data
$ws.Cells.Item(250, 1).Value = 9
$ws.Cells.Item(250, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(250, 3).Value = "Metropolitana"
$ws.Cells.Item(250, 4).Value = 44809
$ws.Cells.Item(250, 5).Value = 13
$ws.Cells.Item(250, 6).Value = 100112021
$ws.Cells.Item(250, 7).Value = "Ají"
$ws.Cells.Item(250, 8).Value = "Inferno"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 60
$ws.Cells.Item(250, 11).Value = 20000
$ws.Cells.Item(250, 12).Value = 20000
$ws.Cells.Item(250, 13).Value = 20000
$ws.Cells.Item(250, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(250, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(250, 16).Value = 2000
$ws.Cells.Item(250, 17).Value = 10
$ws.Cells.Item(250, 18).Value = "Hortaliza"

# New row 251 data
$ws.Cells.Item(251, 1).Value = 9
$ws.Cells.Item(251, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(251, 3).Value = "Metropolitana"
$ws.Cells.Item(251, 4).Value = 44809
$ws.Cells.Item(251, 5).Value = 13
$ws.Cells.Item(251, 6).Value = 100112021
$ws.Cells.Item(251, 7).Value = "Ají"
$ws.Cells.Item(251, 8).Value = "Inferno"
$ws.Cells.Item(251, 9).Value = "Segunda"
$ws.Cells.Item(251, 10).Value = 35
$ws.Cells.Item(251, 11).Value = 17000
$ws.Cells.Item(251, 12).Value = 17000
$ws.Cells.Item(251, 13).Value = 17000
$ws.Cells.Item(251, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(251, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(251, 16).Value = 1700
$ws.Cells.Item(251, 17).Value = 10
$ws.Cells.Item(251, 18).Value = "Hortaliza"
